$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.622.58"
$ws.Range("E2").Value = "  -3.51%  "
$ws.Range("D3").Value = "2.602.93"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'574.07"
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("D6").Value = "'155.53"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -3.79%  "
$ws.Range("E9").Value = "  -7.44%  "
$ws.Range("D10").Value = "'5.86"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "'0.379"
$ws.Range("E11").Value = "  -5.86%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "'28.02"
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("D14").Value = "3.067.72"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "'0.0000178"
$ws.Range("E15").Value = "  -9.52%  "
$ws.Range("D16").Value = "63.444.90"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "2.571.57"
$ws.Range("E17").Value = "  -4.11%  "
$ws.Range("D18").Value = "'12.01"
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").Value = "'7.54"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'4.53"
$ws.Range("E20").Value = "  -6.49%  "
$ws.Range("D21").Value = "'342.19"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D23").Value = "'67.47"
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("E25").Value = "  -5.72%  "
$ws.Range("D26").Value = "'9.15"
$ws.Range("E26").Value = "  -6.06%  "
$ws.Range("D27").Value = "'580.71"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.161"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'7.90"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("D34").Value = "'6.53"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  -4.23%  "
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'19.72"
$ws.Range("E38").Value = "  -4.55%  "
$ws.Range("D39").Value = "'154.40"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "'1.87"
$ws.Range("E40").Value = "  -5.17%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'41.27"
$ws.Range("E42").Value = "  -3.96%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.46"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").Value = "'157.29"
$ws.Range("E44").Value = "  -3.18%  "
$ws.Range("E45").Value = "  -5.67%  "
$ws.Range("D46").Value = "'23.44"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  -5.55%  "
$ws.Range("D48").Value = "'0.627"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").Value = "'18.79"
$ws.Range("E51").Value = "  -5.67%  "
